$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# Row 3 - Title / Romeo & Juliet
$ws.Range("B3").Value = "Title"
$ws.Range("D3").Value = "Romeo & Juliet"

# Row 4 - Author / Shakespeare
$ws.Range("B4").Value = "Author"
$ws.Range("D4").Value = "Shakespeare"

# Row 6 - Num Sold / 300
$ws.Range("B6").Value = "Num Sold "
$ws.Range("D6").Value = 300

# Row 8 - Date Issued / date value formatted as short date
$ws.Range("B8").Value = "Date Issued"
$ws.Range("D8").Value = 39815
$ws.Range("D8").NumberFormat = "mm-dd-yy"

# Row 9 - Date Issued (Error) / text pseudo-date with custom formatted style
$ws.Range("B9").Value = "Date Issued (Error)"
$ws.Range("D9").Value = " 8/18/06"
$ws.Range("D9").Font.Name = "Tahoma"
$ws.Range("D9").Font.Size = 10
$ws.Range("D9").Interior.ColorIndex = 42
$ws.Range("D9").Borders.LineStyle = 1
$ws.Range("D9").NumberFormat = "[$-F800]dddd, mmmm dd, yyyy"
$ws.Range("D9").HorizontalAlignment = -4108
$ws.Range("D9").Locked = $false

# F9 left empty but carries the same date style as D8
$ws.Range("F9").NumberFormat = "mm-dd-yy"

$ws.Range("A1:XFD1048576").Select()
